$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refactor: rename the "trafo_id" column header to "gridnode_id"
$ws.Range("J1").Value = "gridnode_id"

# Move the active selection as reflected in the saved file
$ws.Range("G6").Select()
